# Applies the stratification_info.xlsx data-fix commit:
#  - fills in a missing age value (B140)
#  - corrects several mis-ordered/incorrect age, gender and group values
#    in rows 150-151 and 316-345
#  - removes a duplicated/shifted data row (the old row 339 numbers are
#    dropped and row 347 is deleted entirely, shrinking the table from
#    A1:D347 to A1:D346)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 140 was missing an age value
$ws.Range("B140").Value = 32

# Rows 150/151: ages were swapped
$ws.Range("B150").Value = 26
$ws.Range("B151").Value = 25

# Rows 316-319: age values shifted
$ws.Range("B316").Value = 19
$ws.Range("B317").Value = 27
$ws.Range("B318").Value = 42
$ws.Range("B319").Value = 18

# Row 321: age + group corrected
$ws.Range("B321").Value = 24
$ws.Range("D321").Value = "patient"

# Row 322: age corrected
$ws.Range("B322").Value = 28

# Rows 323-334: ages/genders corrected and missing "group" values filled in
$ws.Range("B323").Value = 30
$ws.Range("D323").Value = "patient"

$ws.Range("B324").Value = 20
$ws.Range("C324").Value = "male"
$ws.Range("D324").Value = "patient"

$ws.Range("B325").Value = 36
$ws.Range("C325").Value = "female"
$ws.Range("D325").Value = "patient"

$ws.Range("B326").Value = 52
$ws.Range("D326").Value = "patient"

$ws.Range("B327").Value = 30
$ws.Range("C327").Value = "male"
$ws.Range("D327").Value = "patient"

$ws.Range("B328").Value = 28
$ws.Range("C328").Value = "female"
$ws.Range("D328").Value = "patient"

$ws.Range("B329").Value = 37
$ws.Range("D329").Value = "patient"

$ws.Range("B330").Value = 26
$ws.Range("D330").Value = "HC"

$ws.Range("B331").Value = 39
$ws.Range("C331").Value = "male"
$ws.Range("D331").Value = "HC"

$ws.Range("B332").Value = 58
$ws.Range("D332").Value = "HC"

$ws.Range("B333").Value = 23
$ws.Range("D333").Value = "patient"

$ws.Range("B334").Value = 36
$ws.Range("D334").Value = "HC"

# Rows 335-337: missing "group" values filled in
$ws.Range("D335").Value = "HC"
$ws.Range("D336").Value = "HC"
$ws.Range("D337").Value = "HC"

# Row 338: project id, age and group corrected/filled in
$ws.Range("A338").Value = "9"
$ws.Range("B338").Value = 32
$ws.Range("D338").Value = "patient"

# Row 339: age/gender removed (data moved to row 340), group stays "patient"
$ws.Range("B339").ClearContents()
$ws.Range("C339").ClearContents()

# Row 340: age/gender filled in
$ws.Range("B340").Value = 25
$ws.Range("C340").Value = "female"

# Rows 341-345: ages/genders corrected
$ws.Range("B341").Value = 59

$ws.Range("B342").Value = 25
$ws.Range("C342").Value = "male"

$ws.Range("B343").Value = 26
$ws.Range("C343").Value = "female"

$ws.Range("B344").Value = 35
$ws.Range("C344").Value = "male"

$ws.Range("B345").Value = 40
$ws.Range("C345").Value = "female"

# Row 346 is unchanged; row 347 is removed entirely (duplicate trailing row)
$ws.Rows(347).Delete()
